$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pluralize "cell line" -> "cell lines" and "iPSC, differentiated" -> "iPSCs, differentiated"
# in the count/type summary column (E) of the milestones table.
$ws.Range("E4").Value = '3\cell lines'
$ws.Range("E5").Value = '3\cell lines'
$ws.Range("E6").Value = '3\cell lines'
$ws.Range("E7").Value = '1\cell lines;2\iPSCs, differentiated'
$ws.Range("E8").Value = '1\cell lines;2\iPSCs, differentiated'
$ws.Range("E9").Value = '1\cell lines;2\iPSCs, differentiated'
$ws.Range("E10").Value = '2\iPSCs, differentiated'
$ws.Range("E11").Value = '2\iPSCs, differentiated'
$ws.Range("E12").Value = '2\iPSCs, differentiated'
$ws.Range("E13").Value = '4\cell lines'
$ws.Range("E14").Value = '4\cell lines'
$ws.Range("E15").Value = '4\cell lines'
$ws.Range("E16").Value = '2\iPSCs, differentiated'
$ws.Range("E17").Value = '2\iPSCs, differentiated'
$ws.Range("E18").Value = '2\iPSCs, differentiated'
$ws.Range("E20").Value = '6\cell lines'
$ws.Range("E21").Value = '6\cell lines'
$ws.Range("E22").Value = '7\cell lines'
$ws.Range("E23").Value = '7\cell lines'
$ws.Range("E25").Value = '2\cell lines'
$ws.Range("E26").Value = '2\cell lines'
$ws.Range("E27").Value = '10\cell lines'
$ws.Range("E28").Value = '10\cell lines'
$ws.Range("E29").Value = '4\cell lines'
$ws.Range("E30").Value = '2\cell lines'
$ws.Range("E31").Value = '2\cell lines'
$ws.Range("E32").Value = '6\cell lines'
$ws.Range("E33").Value = '6\cell lines'
$ws.Range("E35").Value = '6\cell lines'
$ws.Range("E36").Value = '4\cell lines'
$ws.Range("E37").Value = '2\cell lines'
$ws.Range("E38").Value = '2\cell lines'
$ws.Range("E39").Value = '7\cell lines'
$ws.Range("E40").Value = '1\cell lines'
$ws.Range("E41").Value = '5\cell lines'
$ws.Range("E42").Value = '5\cell lines'
$ws.Range("E43").Value = '4\cell lines'
$ws.Range("E44").Value = '6\cell lines'
$ws.Range("E46").Value = '6\cell lines'
$ws.Range("E47").Value = '6\cell lines'
$ws.Range("E49").Value = '12\iPSCs, differentiated'
$ws.Range("E50").Value = '12\iPSCs, differentiated'
$ws.Range("E51").Value = '12\iPSCs, differentiated'
$ws.Range("E52").Value = '39\iPSCs, differentiated'
$ws.Range("E53").Value = '12\iPSCs, differentiated'
$ws.Range("E55").Value = '12\iPSCs, differentiated'
$ws.Range("E56").Value = '12\iPSCs, differentiated;12\mature motor neurons'
$ws.Range("E57").Value = '12\iPSCs, differentiated;12\mature motor neurons'
$ws.Range("E58").Value = '12\iPSCs, differentiated;12\mature motor neurons'
$ws.Range("E59").Value = '12\iPSCs, differentiated'
$ws.Range("E60").Value = '12\iPSCs, differentiated'
$ws.Range("E61").Value = '12\iPSCs, differentiated'
$ws.Range("E63").Value = '3\cell lines'
$ws.Range("E64").Value = '3\cell lines'

# Update the active selection/scroll position to column E (as left by the editor)
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Columns("E:E").Select()
